$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 20 - Aula 31 "O Processo de fragmentação"
# ---------------------------------------------------------------------
$ws.Range("B20").Value = 31
$ws.Range("C20").Value = "7. Thymeleaf para as Views"

# note: E20 is populated before D20 to mirror the shared-string order
# recorded by the authoring tool (the long note became string index 40,
# the short title became string index 41)
$run1 = "nenhum código foi criado, porém explica de forma bem didática e detalhada o funcionamento e a aplicação de fragmentos em paginas HTML. No meu entendimento, fragments são pedaços ou componentes de codigos em documentos HTML que usando as tags de acordo, podem ser reaproveitadas em outras páginas. Existêm 3 formas de reaproveitamento que são: "
$run2 = "INSERT, INCLUDE e REPLACE."
$run3 = " Cada um terá um comportamento diferente na página de destino. INSERT inclui a tag inteira e seu conteudo. REPLACE substitui o componente pai e insere somente o componente do fragmento. REPLACE ignora a tag pai e a tag do fragmento e insere somente o seu VALOR/CONTEUDO."

$ws.Range("E20").Value = $run1 + $run2 + $run3
$ws.Range("E20").WrapText = $true

# highlight the middle run (INSERT, INCLUDE e REPLACE.) in bold red
$startPos = $run1.Length + 1
$len = $run2.Length
$chars = $ws.Range("E20").Characters($startPos, $len)
$chars.Font.Bold = $true
$chars.Font.Color = 255

$ws.Range("D20").Value = "31. O Processo de fragmentação"

$ws.Rows(20).RowHeight = 135

# ---------------------------------------------------------------------
# Row 21 - Aula 32 "Fragmentando o cabeçalho e rodapé"
# ---------------------------------------------------------------------
$ws.Range("B21").Value = 32
$ws.Range("C21").Value = "7. Thymeleaf para as Views"
$ws.Range("D21").Value = "32 . Fragmentando o cabeçalho e rodapé"
$ws.Range("E21").Value = "mostra na prática como é feito o processo de fragmentação, como é feito a chamada entre documentos HTML usando tag do thymeleaf"
$ws.Range("E21").WrapText = $true

$ws.Rows(21).RowHeight = 30

# ---------------------------------------------------------------------
# Update selection to match the newly entered data
# ---------------------------------------------------------------------
[void]$ws.Range("E21").Select()
